$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 10.799999999999986
$ws.Range("B2").Value = 7.0099883166861394
$ws.Range("C2").Value = "aperiodic_7.0.mat"

$ws.Range("A3").Value = 12.499999999999984
$ws.Range("B3").Value = 7.9966533389111021
$ws.Range("C3").Value = "aperiodic_8.0.mat"

$ws.Range("A4").Value = 12.799999999999983
$ws.Range("B4").Value = 8.0933198444669259
$ws.Range("C4").Value = "aperiodic_8.1.mat"

$ws.Range("A5").Value = 12.799999999999983
$ws.Range("B5").Value = 8.1933196778005382
$ws.Range("C5").Value = "aperiodic_8.2.mat"

$ws.Range("A6").Value = 12.899999999999984
$ws.Range("B6").Value = 8.2933195111341487
$ws.Range("C6").Value = "aperiodic_8.3.mat"

$ws.Range("A7").Value = 13.199999999999983
$ws.Range("B7").Value = 8.4066526555789078
$ws.Range("C7").Value = "aperiodic_8.4.mat"

$ws.Range("A8").Value = 13.299999999999983
$ws.Range("B8").Value = 8.5033191611347316
$ws.Range("C8").Value = "aperiodic_8.5.mat"

$ws.Range("A9").Value = 13.599999999999982
$ws.Range("B9").Value = 8.5966523389127687
$ws.Range("C9").Value = "aperiodic_8.6.mat"

$ws.Range("A10").Value = 13.299999999999983
$ws.Range("B10").Value = 8.6999855000241677
$ws.Range("C10").Value = "aperiodic_8.7.mat"

$ws.Range("A11").Value = 13.899999999999983
$ws.Range("B11").Value = 8.8066519889133517
$ws.Range("C11").Value = "aperiodic_8.8.mat"

$ws.Range("A12").Value = 13.999999999999982
$ws.Range("B12").Value = 8.9099851500247507
$ws.Range("C12").Value = "aperiodic_8.9.mat"

$ws.Range("A13").Value = 13.999999999999982
$ws.Range("B13").Value = 8.9966516722472125
$ws.Range("C13").Value = "aperiodic_9.0.mat"

$ws.Range("A14").Value = 14.299999999999981
$ws.Range("B14").Value = 9.1033181611363982
$ws.Range("C14").Value = "aperiodic_9.1.mat"

$ws.Range("A15").Value = 14.399999999999981
$ws.Range("B15").Value = 9.1966513389144353
$ws.Range("C15").Value = "aperiodic_9.2.mat"

$ws.Range("A16").Value = 14.799999999999981
$ws.Range("B16").Value = 9.2999845000258343
$ws.Range("C16").Value = "aperiodic_9.3.mat"

$ws.Range("A17").Value = 14.799999999999981
$ws.Range("B17").Value = 9.3999843333594448
$ws.Range("C17").Value = "aperiodic_9.4.mat"

$ws.Range("A18").Value = 14.99999999999998
$ws.Range("B18").Value = 9.4999841666930553
$ws.Range("C18").Value = "aperiodic_9.5.mat"

$ws.Range("A19").Value = 15.299999999999979
$ws.Range("B19").Value = 9.6033173278044544
$ws.Range("C19").Value = "aperiodic_9.6.mat"

$ws.Range("A20").Value = 15.299999999999979
$ws.Range("B20").Value = 9.6933171778047047
$ws.Range("C20").Value = "aperiodic_9.7.mat"

$ws.Range("A21").Value = 15.099999999999982
$ws.Range("B21").Value = 9.7999836666938887
$ws.Range("C21").Value = "aperiodic_9.8.mat"

$ws.Range("A22").Value = 15.59999999999998
$ws.Range("B22").Value = 9.9099834833608611
$ws.Range("C22").Value = "aperiodic_9.9.mat"

$ws.Range("A23").Value = 15.799999999999979
$ws.Range("B23").Value = 10.003316661138898
$ws.Range("C23").Value = "aperiodic_10.0.mat"

$ws.Range("A24").Value = 15.799999999999979
$ws.Range("B24").Value = 10.10331649447251
$ws.Range("C24").Value = "aperiodic_10.1.mat"

$ws.Range("A25").Value = 16.299999999999983
$ws.Range("B25").Value = 10.206649655583908
$ws.Range("C25").Value = "aperiodic_10.2.mat"

$ws.Range("A26").Value = 16.299999999999983
$ws.Range("B26").Value = 10.303316161139731
$ws.Range("C26").Value = "aperiodic_10.3.mat"

$ws.Range("A27").Value = 16.499999999999986
$ws.Range("B27").Value = 10.409982650028917
$ws.Range("C27").Value = "aperiodic_10.4.mat"

$ws.Range("A28").Value = 16.599999999999987
$ws.Range("B28").Value = 10.493315844473592
$ws.Range("C28").Value = "aperiodic_10.5.mat"

$ws.Range("A29").Value = 16.79999999999999
$ws.Range("B29").Value = 10.593315677807205
$ws.Range("C29").Value = "aperiodic_10.6.mat"

$ws.Range("A30").Value = 16.699999999999989
$ws.Range("B30").Value = 10.693315511140815
$ws.Range("C30").Value = "aperiodic_10.7.mat"

$ws.Range("A31").Value = 17.300000000000001
$ws.Range("B31").Value = 10.809981983363361
$ws.Range("C31").Value = "aperiodic_10.8.mat"

$ws.Range("A32").Value = 17.199999999999999
$ws.Range("B32").Value = 10.903315161141398
$ws.Range("C32").Value = "aperiodic_10.9.mat"

$ws.Range("A33").Value = 17.300000000000001
$ws.Range("B33").Value = 10.996648338919435
$ws.Range("C33").Value = "aperiodic_11.0.mat"

$ws.Range("A34").Value = 17.600000000000005
$ws.Range("B34").Value = 11.109981483364194
$ws.Range("C34").Value = "aperiodic_11.1.mat"

$ws.Range("A35").Value = 17.800000000000008
$ws.Range("B35").Value = 11.203314661142231
$ws.Range("C35").Value = "aperiodic_11.2.mat"

$ws.Range("A36").Value = 17.900000000000009
$ws.Range("B36").Value = 11.293314511142482
$ws.Range("C36").Value = "aperiodic_11.3.mat"

$ws.Range("A37").Value = 18.000000000000011
$ws.Range("B37").Value = 11.409980983365028
$ws.Range("C37").Value = "aperiodic_11.4.mat"

$ws.Range("A38").Value = 18.300000000000015
$ws.Range("B38").Value = 11.506647488920851
$ws.Range("C38").Value = "aperiodic_11.5.mat"

$ws.Range("A39").Value = 18.400000000000016
$ws.Range("B39").Value = 11.606647322254464
$ws.Range("C39").Value = "aperiodic_11.6.mat"

$ws.Range("A40").Value = 18.800000000000022
$ws.Range("B40").Value = 11.706647155588074
$ws.Range("C40").Value = "aperiodic_11.7.mat"

$ws.Range("A41").Value = 18.600000000000023
$ws.Range("B41").Value = 11.803313661143898
$ws.Range("C41").Value = "aperiodic_11.8.mat"

$ws.Range("A42").Value = 19.000000000000028
$ws.Range("B42").Value = 11.893313511144148
$ws.Range("C42").Value = "aperiodic_11.9.mat"

$ws.Range("A43").Value = 19.10000000000003
$ws.Range("B43").Value = 12.003313327811121
$ws.Range("C43").Value = "aperiodic_12.0.mat"

$ws.Range("A44").Value = 24.500000000000117
$ws.Range("B44").Value = 14.99664167226388
$ws.Range("C44").Value = "aperiodic_15.0.mat"

